# Update Log Book dan Burndown Chart
# Fill in the "2/6", "3/6", "4/6" (columns L, M, N) remaining-effort values
# for each task row (5-14) on Sheet1, then move the active selection to
# reflect where the user was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0

$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0

$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0

$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 0

$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 0

$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0

$ws.Range("L11").Value = 3
$ws.Range("M11").Value = 2
$ws.Range("N11").Value = 0

$ws.Range("L12").Value = 4
$ws.Range("M12").Value = 3
$ws.Range("N12").Value = 0

$ws.Range("L13").Value = 2
$ws.Range("M13").Value = 2
$ws.Range("N13").Value = 2

$ws.Range("L14").Value = 5
$ws.Range("M14").Value = 5
$ws.Range("N14").Value = 4

$ws.Activate()
$ws.Range("N14").Select()
